$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_AE")
$ws.Name = "CRF_AE"
